$d = $word.ActiveDocument

# --- 1. "Kawarazuka" run: remove underline (add <w:u w:val="none"/>) ---
$rng = $d.Content
$rng.Find.Execute("Kawarazuka", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Underline = 0

# --- 2. " and " run: remove underline ---
$rng = $d.Content
$rng.Find.Execute(" and ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Underline = 0

# --- 3. "Béné" run: remove underline ---
$rng = $d.Content
$rng.Find.Execute("Béné", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Underline = 0

# --- 4. " 2010" run: text -> ", 2010" and remove underline ---
$rng = $d.Content
$rng.Find.Execute(" 2010", $true, $false, $false, $false, $false, $true, 1, $false, ", 2010", 2)
$rng = $d.Content
$rng.Find.Execute(", 2010", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Underline = 0

# --- 5. "Allison 2011" -> "Allison, 2011" and remove underline ---
$rng = $d.Content
$rng.Find.Execute("Allison 2011", $true, $false, $false, $false, $false, $true, 1, $false, "Allison, 2011", 2)
$rng = $d.Content
$rng.Find.Execute("Allison, 2011", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Underline = 0

# --- 6. "Golden et al. 2016" -> "Golden et al., 2016" and remove underline ---
$rng = $d.Content
$rng.Find.Execute("Golden et al. 2016", $true, $false, $false, $false, $false, $true, 1, $false, "Golden et al., 2016", 2)
$rng = $d.Content
$rng.Find.Execute("Golden et al., 2016", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Underline = 0

Write-Output $d.Content.Text
